$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6891
$ws.Range("J3").Value = 7292
$ws.Range("C4").Value = 1842
$ws.Range("F4").Value = 1903
$ws.Range("I4").Value = 1775
$ws.Range("J4").Value = 1585
$ws.Range("J5").Value = 575
$ws.Range("J6").Value = 9762
$ws.Range("C7").Value = 28386
$ws.Range("F7").Value = 24094
$ws.Range("I7").Value = 26232
$ws.Range("J7").Value = 26105

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 252
$ws.Range("J7").Value = 380

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 439
$ws.Range("J3").Value = 488
$ws.Range("J6").Value = 584
$ws.Range("J7").Value = 1641

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 139
$ws.Range("J7").Value = 523

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 392
$ws.Range("J5").Value = 50
$ws.Range("J7").Value = 1182

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 133
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 381

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 241
$ws.Range("J3").Value = 272
$ws.Range("J4").Value = 33
$ws.Range("J6").Value = 234
$ws.Range("J7").Value = 810

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 188
$ws.Range("J6").Value = 239
$ws.Range("J7").Value = 654

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 206
$ws.Range("J7").Value = 755
$ws.Range("J8").Value = 1641
$ws.Range("J11").Value = 454
$ws.Range("J15").Value = 311
$ws.Range("J19").Value = 760
$ws.Range("J20").Value = 547
$ws.Range("J27").Value = 154
$ws.Range("J29").Value = 1411
$ws.Range("J30").Value = 92
$ws.Range("J31").Value = 263
$ws.Range("J33").Value = 1182
$ws.Range("J34").Value = 115
$ws.Range("J36").Value = 357
$ws.Range("J37").Value = 810
$ws.Range("J42").Value = 1126
$ws.Range("J43").Value = 224
$ws.Range("J48").Value = 297
$ws.Range("J52").Value = 664
$ws.Range("J53").Value = 380
$ws.Range("J55").Value = 406
$ws.Range("J60").Value = 152
$ws.Range("C63").Value = 271
$ws.Range("F63").Value = 188
$ws.Range("I63").Value = 183
$ws.Range("J65").Value = 654
$ws.Range("J66").Value = 80
$ws.Range("J67").Value = 975
$ws.Range("J72").Value = 102
$ws.Range("J78").Value = 308
$ws.Range("J79").Value = 732
$ws.Range("J83").Value = 523
$ws.Range("J84").Value = 218
$ws.Range("J85").Value = 1073
$ws.Range("J88").Value = 276
$ws.Range("J90").Value = 278
$ws.Range("J94").Value = 281
$ws.Range("J95").Value = 381
$ws.Range("C101").Value = 28386
$ws.Range("F101").Value = 24094
$ws.Range("I101").Value = 26232
$ws.Range("J101").Value = 26105

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 362
$ws.Range("J7").Value = 975

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 66
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 501
$ws.Range("J4").Value = 74
$ws.Range("J7").Value = 1411

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 183
$ws.Range("J3").Value = 217
$ws.Range("J6").Value = 294
$ws.Range("J7").Value = 760

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 220
$ws.Range("J6").Value = 600
$ws.Range("J7").Value = 1126

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 97
$ws.Range("J7").Value = 308

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J5").Value = 6
$ws.Range("J7").Value = 406

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J6").Value = 220
$ws.Range("J7").Value = 732

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J4").Value = 43
$ws.Range("J5").Value = 12
$ws.Range("J7").Value = 547

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 114
$ws.Range("J7").Value = 357

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 236
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 755

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 54
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 89
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J4").Value = 27
$ws.Range("J6").Value = 210
$ws.Range("J7").Value = 454

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J3").Value = 50
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 140
$ws.Range("J7").Value = 276

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 36
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 99
$ws.Range("J7").Value = 278

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 53
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 382
$ws.Range("J7").Value = 1073

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 102

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 283
$ws.Range("J7").Value = 664
